$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New query text for the FilesTab stats query (row 3, column C)
$filesStatQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
   WHERE c.race = "ASIAN"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

# New query text for the FilesTab data query (row 3, column B)
$filesDataQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
     WHERE c.race = "ASIAN"
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

# New query text for the CasesTab data query (row 2, column B)
$casesDataQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
   WHERE c.race = "ASIAN"
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

# New query text for the CasesTab stats query (row 2, column C)
$casesStatQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
      WHERE c.race = "ASIAN"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

# Add row 3 (FilesTab) first, in the same order the new shared strings were
# originally authored in the workbook (A3, C3, B3, then the row 2 updates).
$ws.Cells.Item(3, 1).Value = "FilesTab"
$ws.Cells.Item(3, 3).Value = $filesStatQuery
$ws.Cells.Item(3, 2).Value = $filesDataQuery
$ws.Cells.Item(3, 4).Value = $ws.Cells.Item(2, 4).Value()
$ws.Cells.Item(3, 5).Value = $ws.Cells.Item(2, 5).Value()

# Update row 2 (CasesTab) to use the new single-filter-script queries
$ws.Cells.Item(2, 2).Value = $casesDataQuery
$ws.Cells.Item(2, 3).Value = $casesStatQuery

$ws.Range("B2:C3").WrapText = $true

# Row heights grow to fit the longer wrapped query text (row 3 hits Excel's
# maximum row height of 409.5pt).
$ws.Rows.Item(2).RowHeight = 195
$ws.Rows.Item(3).RowHeight = 409.5

$ws.Range("B3").Select() | Out-Null
